# Region VI_HEALTH.xlsx edit:
# 1. Rename header B1 from "REGION" to "Region"
# 2. Narrow column C width from 19 to 17
# 3. Remove data row 11 (Negros Occidental / Efegenio Lizares NHS), shifting
#    every subsequent row up by one (dimension/validation ranges follow
#    automatically via the row delete).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Header rename
$ws.Range("B1").Value = "Region"

# 2. Column width change (column C = index 3).
#    Excel's ColumnWidth (character units) maps to the raw OOXML <col width>
#    with a constant +5/6 padding offset in this engine, so back that out to
#    land exactly on width="17" in the saved file.
$ws.Columns.Item(3).ColumnWidth = (17 - (5/6))

# 3. Delete entire row 11
$ws.Rows.Item(11).Delete()
